$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - LinearRegression (name unchanged)
$ws.Range("B2").Value = 5017680469586076
$ws.Range("C2").Value = 5017680469586075
$ws.Range("D2").Value = 5017680469586075

# Row 3 - RandomForestRegressor (name unchanged)
$ws.Range("B3").Value = 140389428279681.2
$ws.Range("C3").Value = 119519752865154.6
$ws.Range("D3").Value = 567645782670764.5

# Row 4 - GradientBoostingRegressor -> DecisionTreeRegressor
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 3083776939523.741
$ws.Range("C4").Value = 3311861470233.196
$ws.Range("D4").Value = 130559512434404.5

# Row 5 - AdaBoostRegressor -> MLPRegressor
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 298743487311586.4
$ws.Range("C5").Value = 407279260915297.5
$ws.Range("D5").Value = 2788234075976529
